$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# SimulationParameters: add "nra" renamed row stays same content, but add two
# new rows (Meshwidth / h) and (Change when LOS, PerfRef or Tx are changed /
# Runplottype / PerfectRefCentre) below the existing data.
# ---------------------------------------------------------------------------
$wsSim = $wb.Worksheets.Item("SimulationParameters")

$wsSim.Range("A13").Value = "Meshwidth"
$wsSim.Range("B13").Value = "h"

$wsSim.Range("A14").Value = "Change when LOS, PerfRef or Tx are changed"
$wsSim.Range("B14").Value = "Runplottype"
$wsSim.Range("C14").Value = "PerfectRefCentre"

# ---------------------------------------------------------------------------
# Obstacles: drop the merged "Box Objects" banner row, add "Box"/"Triangle"
# columns and the per-vertex triangle-mesh columns (p0x..p2z).
# ---------------------------------------------------------------------------
$wsObs = $wb.Worksheets.Item("Obstacles")

$wsObs.Rows("1:1").Delete()

$wsObs.Range("C1:D1").EntireColumn.Insert()
$wsObs.Range("C1").Value = "Box"
$wsObs.Range("D1").Value = "Triangle"
$wsObs.Range("C2").Value = 1
$wsObs.Range("D2").Value = 0

$wsObs.Range("K1").Value = "p0x"
$wsObs.Range("L1").Value = "p0y"
$wsObs.Range("M1").Value = "p0z"
$wsObs.Range("N1").Value = "p1x"
$wsObs.Range("O1").Value = "p1y"
$wsObs.Range("P1").Value = "p1z"
$wsObs.Range("Q1").Value = "p2x"
$wsObs.Range("R1").Value = "p2y"
$wsObs.Range("S1").Value = "p2z"

# ---------------------------------------------------------------------------
# OuterBoundary: same restructuring as Obstacles.
# ---------------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("OuterBoundary")

$wsOut.Rows("1:1").Delete()

$wsOut.Range("C1:D1").EntireColumn.Insert()
$wsOut.Range("C1").Value = "Box"
$wsOut.Range("D1").Value = "Triangle"
$wsOut.Range("C2").Value = 1
$wsOut.Range("D2").Value = 0

$wsOut.Range("K1").Value = "p0x"
$wsOut.Range("L1").Value = "p0y"
$wsOut.Range("M1").Value = "p0z"
$wsOut.Range("N1").Value = "p1x"
$wsOut.Range("O1").Value = "p1y"
$wsOut.Range("P1").Value = "p1z"
$wsOut.Range("Q1").Value = "p2x"
$wsOut.Range("R1").Value = "p2y"
$wsOut.Range("S1").Value = "p2z"

# ---------------------------------------------------------------------------
# Leave SimulationParameters as the active tab/sheet (matches new activeTab).
# ---------------------------------------------------------------------------
$wsSim.Activate()
$wsSim.Range("A15").Select()
